$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "PackageID" / "Null" data row beneath the existing table
$ws.Range("B17").Value = "PackageID"
$ws.Range("C17").Value = "Null"

# Widen column B to fit the new content (matches the author's manual resize)
$ws.Columns.Item(2).ColumnWidth = 15.45

# Leave the selection where it ends up after typing the last entry
$ws.Range("D17").Select()
